$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-16 Friday" "2026-01-17 Saturday"

Replace-Text "850÷2=425, 0" "463÷8=57, 7"
Replace-Text "984÷9=109, 3" "759÷3=253, 0"
Replace-Text "392÷5=78, 2" "886÷8=110, 6"
Replace-Text "931÷8=116, 3" "120÷5=24, 0"
Replace-Text "648÷6=108, 0" "606÷7=86, 4"

Replace-Text "792÷9=88, 0" "965÷8=120, 5"
Replace-Text "915÷7=130, 5" "977÷9=108, 5"
Replace-Text "271÷4=67, 3" "964÷5=192, 4"
Replace-Text "161÷9=17, 8" "449÷5=89, 4"
Replace-Text "869÷6=144, 5" "378÷3=126, 0"

Replace-Text "788÷5=157, 3" "311÷8=38, 7"
Replace-Text "825÷3=275, 0" "729÷7=104, 1"
Replace-Text "362÷3=120, 2" "108÷3=36, 0"
Replace-Text "675÷4=168, 3" "746÷8=93, 2"
Replace-Text "573÷4=143, 1" "256÷6=42, 4"

Replace-Text "310÷4=77, 2" "139÷7=19, 6"
Replace-Text "660÷9=73, 3" "321÷4=80, 1"
Replace-Text "125÷7=17, 6" "297÷9=33, 0"
Replace-Text "623÷6=103, 5" "995÷8=124, 3"
Replace-Text "964÷3=321, 1" "554÷9=61, 5"

Replace-Text "265÷4=66, 1" "928÷8=116, 0"
Replace-Text "433÷5=86, 3" "743÷9=82, 5"
Replace-Text "472÷5=94, 2" "330÷2=165, 0"
Replace-Text "527÷7=75, 2" "193÷4=48, 1"
Replace-Text "793÷3=264, 1" "833÷3=277, 2"

Write-Output "Done"
